$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").Value = 25
$ws.Range("AB2").Value = 14.5
$ws.Range("AD2").Value = 11
$ws.Range("AE2").Value = 22
$ws.Range("AF2").Value = 28
$ws.Range("AG2").Value = 15.5
$ws.Range("AH2").Value = 19
$ws.Range("AI2").Value = 36
$ws.Range("AJ2").Value = 80
$ws.Range("AK2").Value = 48
$ws.Range("AM2").Value = 90
$ws.Range("AN2").Value = 46
$ws.Range("AO2").Value = 15.5
$ws.Range("F2").Value = 3.85
$ws.Range("G2").Value = 3.95
$ws.Range("H2").Value = 2.1
$ws.Range("I2").Value = 2.14
$ws.Range("J2").Value = 3.7
$ws.Range("K2").Value = 3.75
$ws.Range("L2").Value = 1.41
$ws.Range("N2").Value = 3.95
$ws.Range("O2").Value = 1.32
$ws.Range("P2").Value = 1.98
$ws.Range("Q2").Value = 1.98
$ws.Range("R2").Value = 1.39
$ws.Range("S2").Value = 3.45
$ws.Range("T2").Value = 1.8
$ws.Range("U2").Value = 2.18
$ws.Range("V2").Value = 1.88
$ws.Range("X2").Value = 15.5
$ws.Range("Y2").Value = 9.800000000000001
$ws.Range("Z2").Value = 13
$ws.Range("AC3").Value = 18
$ws.Range("AD3").Value = 60
$ws.Range("F3").Value = 1.25
$ws.Range("H3").Value = 13
$ws.Range("I3").Value = 16.5
$ws.Range("K3").Value = 7.2
$ws.Range("O3").Value = 1.19
$ws.Range("P3").Value = 2.42
$ws.Range("R3").Value = 1.56
$ws.Range("T3").Value = 2.1
$ws.Range("U3").Value = 1.73
$ws.Range("AA4").Value = 48
$ws.Range("AC4").Value = 14.5
$ws.Range("AD4").Value = 20
$ws.Range("AE4").Value = 65
$ws.Range("AH4").Value = 25
$ws.Range("AO4").Value = 55
$ws.Range("G4").Value = 5.9
$ws.Range("H4").Value = 1.7
$ws.Range("I4").Value = 1.83
$ws.Range("K4").Value = 4.5
$ws.Range("N4").Value = 4.8
$ws.Range("P4").Value = 2.3
$ws.Range("Q4").Value = 1.6
$ws.Range("R4").Value = 1.52
$ws.Range("S4").Value = 2.48
$ws.Range("V4").Value = 2.2
$ws.Range("W4").Value = 1.22
$ws.Range("X4").Value = 42
$ws.Range("Y4").Value = 22
$ws.Range("AA5").Value = 25
$ws.Range("AB5").Value = 16.5
$ws.Range("AC5").Value = 8.6
$ws.Range("AD5").Value = 11
$ws.Range("AE5").Value = 23
$ws.Range("AF5").Value = 36
$ws.Range("AG5").Value = 18
$ws.Range("AH5").Value = 19
$ws.Range("AI5").Value = 130
$ws.Range("AJ5").Value = 900
$ws.Range("AK5").Value = 150
$ws.Range("AL5").Value = 330
$ws.Range("AM5").Value = 580
$ws.Range("AO5").Value = 15.5
$ws.Range("F5").Value = 3.95
$ws.Range("H5").Value = 1.92
$ws.Range("I5").Value = 2.04
$ws.Range("J5").Value = 3.5
$ws.Range("K5").Value = 3.95
$ws.Range("L5").Value = 1.34
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 3.6
$ws.Range("O5").Value = 1.3
$ws.Range("P5").Value = 1.88
$ws.Range("R5").Value = 1.34
$ws.Range("U5").Value = 2.04
$ws.Range("V5").Value = 1.96
$ws.Range("W5").Value = 1.27
$ws.Range("X5").Value = 15.5
$ws.Range("Y5").Value = 9.6
$ws.Range("Z5").Value = 13
$ws.Range("AC6").Value = 17.5
$ws.Range("AH6").Value = 60
$ws.Range("AO6").Value = 85
$ws.Range("F6").Value = 3.4
$ws.Range("G6").Value = 4.6
$ws.Range("H6").Value = 1.95
$ws.Range("I6").Value = 2.2
$ws.Range("K6").Value = 4.7
$ws.Range("N6").Value = 3.95
$ws.Range("P6").Value = 2.3
$ws.Range("R6").Value = 1.58
$ws.Range("U6").Value = 2.36
$ws.Range("V6").Value = 1.83
$ws.Range("W6").Value = 1.28
$ws.Range("AA7").Value = 900
$ws.Range("AC7").Value = 22
$ws.Range("AD7").Value = 20
$ws.Range("AE7").Value = 65
$ws.Range("AO7").Value = 15
$ws.Range("F7").Value = 7.2
$ws.Range("G7").Value = 9.6
$ws.Range("H7").Value = 1.46
$ws.Range("I7").Value = 1.48
$ws.Range("J7").Value = 4.7
$ws.Range("K7").Value = 5.5
$ws.Range("N7").Value = 4.7
$ws.Range("P7").Value = 2.28
$ws.Range("Q7").Value = 1.68
$ws.Range("R7").Value = 1.47
$ws.Range("S7").Value = 2.72
$ws.Range("T7").Value = 1.91
$ws.Range("U7").Value = 1.9
$ws.Range("V7").Value = 3.05
$ws.Range("W7").Value = 1.12
$ws.Range("X7").Value = 32
$ws.Range("Y7").Value = 990
$ws.Range("Z7").Value = 17.5
$ws.Range("AA8").Value = 27
$ws.Range("AJ8").Value = 60
$ws.Range("AN8").Value = 21
$ws.Range("F8").Value = 3.45
$ws.Range("H8").Value = 2.04
$ws.Range("I8").Value = 2.14
$ws.Range("K8").Value = 4.5
$ws.Range("P8").Value = 2.66
$ws.Range("R8").Value = 1.67
$ws.Range("S8").Value = 2.24
$ws.Range("V8").Value = 1.89
$ws.Range("W8").Value = 1.39
$ws.Range("Z8").Value = 18
$ws.Range("AA9").Value = 900
$ws.Range("AB9").Value = 14
$ws.Range("AD9").Value = 17.5
$ws.Range("AE9").Value = 95
$ws.Range("AF9").Value = 16
$ws.Range("AG9").Value = 11.5
$ws.Range("AH9").Value = 16.5
$ws.Range("AI9").Value = 130
$ws.Range("AK9").Value = 19.5
$ws.Range("AL9").Value = 60
$ws.Range("AM9").Value = 200
$ws.Range("AO9").Value = 34
$ws.Range("F9").Value = 1.96
$ws.Range("G9").Value = 2.08
$ws.Range("J9").Value = 3.9
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 5.3
$ws.Range("P9").Value = 2.66
$ws.Range("Q9").Value = 1.56
$ws.Range("U9").Value = 2.5
$ws.Range("W9").Value = 1.92
$ws.Range("Y9").Value = 22
$ws.Range("Z9").Value = 80
$ws.Range("AA10").Value = 30
$ws.Range("AB10").Value = 27
$ws.Range("AC10").Value = 9.4
$ws.Range("AD10").Value = 11
$ws.Range("AE10").Value = 19.5
$ws.Range("AF10").Value = 980
$ws.Range("AG10").Value = 42
$ws.Range("AH10").Value = 38
$ws.Range("AI10").Value = 80
$ws.Range("AJ10").Value = 900
$ws.Range("AK10").Value = 190
$ws.Range("AL10").Value = 190
$ws.Range("AM10").Value = 120
$ws.Range("F10").Value = 5.3
$ws.Range("G10").Value = 6.2
$ws.Range("H10").Value = 1.66
$ws.Range("I10").Value = 1.72
$ws.Range("K10").Value = 4.4
$ws.Range("L10").Value = 1.35
$ws.Range("N10").Value = 4.5
$ws.Range("O10").Value = 1.27
$ws.Range("P10").Value = 2.04
$ws.Range("Q10").Value = 1.83
$ws.Range("R10").Value = 1.39
$ws.Range("S10").Value = 3.05
$ws.Range("T10").Value = 1.84
$ws.Range("U10").Value = 2.02
$ws.Range("V10").Value = 2.38
$ws.Range("W10").Value = 1.2
$ws.Range("X10").Value = 30
$ws.Range("Y10").Value = 17.5
$ws.Range("AA11").Value = 11.5
$ws.Range("AC11").Value = 27
$ws.Range("AF11").Value = 170
$ws.Range("AG11").Value = 60
$ws.Range("AH11").Value = 34
$ws.Range("AI11").Value = 36
$ws.Range("AL11").Value = 470
$ws.Range("AM11").Value = 320
$ws.Range("AO11").Value = 3.1
$ws.Range("H11").Value = 1.22
$ws.Range("I11").Value = 1.24
$ws.Range("J11").Value = 7.4
$ws.Range("K11").Value = 8.800000000000001
$ws.Range("R11").Value = 2.02
$ws.Range("S11").Value = 1.8
$ws.Range("U11").Value = 2.06
$ws.Range("V11").Value = 5
$ws.Range("W11").Value = 1.07
$ws.Range("Z11").Value = 11.5
$ws.Range("AA12").Value = 23
$ws.Range("AB12").Value = 26
$ws.Range("AC12").Value = 11
$ws.Range("AD12").Value = 10.5
$ws.Range("AG12").Value = 21
$ws.Range("AK12").Value = 60
$ws.Range("AL12").Value = 60
$ws.Range("AN12").Value = 55
$ws.Range("F12").Value = 4.8
$ws.Range("I12").Value = 1.71
$ws.Range("J12").Value = 4.3
$ws.Range("K12").Value = 4.8
$ws.Range("P12").Value = 2.32
$ws.Range("Q12").Value = 1.64
$ws.Range("S12").Value = 2.56
$ws.Range("V12").Value = 2.4
$ws.Range("W12").Value = 1.23
$ws.Range("Y12").Value = 11.5
$ws.Range("AB13").Value = 12
$ws.Range("AC13").Value = 23
$ws.Range("AF13").Value = 9.199999999999999
$ws.Range("AG13").Value = 11.5
$ws.Range("AH13").Value = 80
$ws.Range("AJ13").Value = 10.5
$ws.Range("AK13").Value = 26
$ws.Range("F13").Value = 1.26
$ws.Range("G13").Value = 1.31
$ws.Range("H13").Value = 12
$ws.Range("I13").Value = 14.5
$ws.Range("J13").Value = 6.2
$ws.Range("K13").Value = 7.4
$ws.Range("M13").Value = 1.02
$ws.Range("N13").Value = 6
$ws.Range("O13").Value = 1.15
$ws.Range("P13").Value = 2.68
$ws.Range("Q13").Value = 1.54
$ws.Range("R13").Value = 1.65
$ws.Range("S13").Value = 2.24
$ws.Range("T13").Value = 1.92
$ws.Range("U13").Value = 1.87
$ws.Range("V13").Value = 1.07
$ws.Range("W13").Value = 4.1
$ws.Range("Y13").Value = 260
